$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell B1 (2 -> 4) and apply bold red-on-yellow formatting
$ws.Range("B1").Value = 4
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Font.Color = 255
$ws.Range("B1").Interior.Color = 65535

# New cell C2 = 3
$ws.Range("C2").Value = 3

# Move the active selection to D9 (matches the saved sheet view)
$null = $ws.Range("D9").Select()
